$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "jogo" rows to append below the existing data (rows 244-253)
$rows = @(
    @{ name = "jogo243"; nums = @(3, 4, 12, 28, 36, 38, 58) },
    @{ name = "jogo244"; nums = @(4, 6, 12, 27, 31, 42, 49) },
    @{ name = "jogo245"; nums = @(1, 7, 27, 38, 48, 52, 57) },
    @{ name = "jogo246"; nums = @(6, 14, 24, 30, 37, 54, 58) },
    @{ name = "jogo247"; nums = @(15, 24, 32, 35, 41, 45, 53) },
    @{ name = "jogo248"; nums = @(4, 18, 20, 26, 40, 53, 58) },
    @{ name = "jogo249"; nums = @(3, 14, 29, 41, 50, 52, 54) },
    @{ name = "jogo250"; nums = @(6, 7, 24, 25, 30, 34, 42) },
    @{ name = "jogo251"; nums = @(22, 26, 27, 44, 47, 49, 59) },
    @{ name = "jogo252"; nums = @(13, 14, 22, 33, 40, 49, 50) }
)

$startRow = 244
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $entry = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $entry.name
    $ws.Cells.Item($r, 1).Font.Underline = $true

    for ($c = 0; $c -lt $entry.nums.Count; $c++) {
        $ws.Cells.Item($r, 2 + $c).Value = $entry.nums[$c]
    }
}

$ws.Application.ActiveWindow.ScrollRow = 226
$ws.Range("R1").Select()
